# "Importar tarjetas terminado (aparentemente)"
#
# Row 10 (ALEX OTALORA's card, cedula 1055314236) is marked as returned to
# DAVID CHAPARRO (cedula 1072642921) on 7/9/2024 — the same "return" pattern
# already recorded for rows 9 and 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Card holder / receiver details for the return leg of row 10.
$ws.Range("I10").Value = 1055314236
$ws.Range("J10").Value = "ALEX OTALORA"
$ws.Range("K10").Value = 1072642921
$ws.Range("L10").Value = "DAVID CHAPARRO"

# M10 was an empty, borderless date cell; once filled it should carry the
# same bordered short-date style already used by the analogous M9/M11
# cells, so copy that formatting over before writing the date value.
$ws.Range("M9").Copy()
$ws.Range("M10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("M10").Value = 45482        # 7/9/2024

# Leave the sheet scrolled to / selecting the row just filled in.
$ws.Activate()
$ws.Range("I10:M10").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2

$wb.Save()
